$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("UCL_CI_Mods(6 lasers)")

# ---------------------------------------------------------------------------
# 1. Populate the previously-blank row 10 with a new "Cables" section header
#    (same look as the other section headers on this sheet, e.g. row 2 / 11).
# ---------------------------------------------------------------------------
$ws3.Range("A10").Value = "Cables"
$ws3.Range("A10").Font.Name = "Calibri"
$ws3.Range("A10").Font.Size = 11
$ws3.Range("A10").Font.Bold = $true
$ws3.Range("A10").Font.Color = RGB(0, 0, 0)
$ws3.Range("A10").Interior.Color = RGB(242, 242, 242)

# ---------------------------------------------------------------------------
# 2. Insert 9 blank rows before (old) row 11 to make room for the new cable
#    parts list (new rows 11-19).
# ---------------------------------------------------------------------------
$ws3.Rows("11:19").Insert()

# Row 11 - Farnell BNC double female adaptor
$ws3.Range("A11").Value = 3384436
$ws3.Range("B11").Value = "Farnell"
$ws3.Range("C11").Value = "BNC DOUBLE FEMALE ADAPTOR ZINC"
$ws3.Range("D11").Value = 14
$ws3.Range("E11").Value = 0.5
$ws3.Range("F11").Formula = "=E11*D11"
$ws3.Range("G11").Value = "Signal generation"

# Row 12 - Farnell BNC plug to free end cable
$ws3.Range("A12").Value = 3703635
$ws3.Range("B12").Value = "Farnell"
$ws3.Range("C12").Value = "BNC PLUG TO FREE END, BLACK/RED, 150MM"
$ws3.Range("D12").Value = 14
$ws3.Range("E12").Value = 6
$ws3.Range("F12").Formula = "=E12*D12"
$ws3.Range("G12").Value = "Signal generation"
$ws3.Range("H12").Value = "male to free end. You could buy female to free end and not buy the adaptors(above) but this was strangely much harder to find and more expensive. "

# Row 13 - Thorlabs SMB coaxial cable
$ws3.Range("A13").Value = "PAA236R"
$ws3.Range("B13").Value = "Thorlabs"
$ws3.Range("C13").Value = "SMB Coaxial Cable, 90° SMB Female to BNC Male, 36"" (914 mm)"
$ws3.Range("D13").Value = 12
$ws3.Range("E13").Value = 12
$ws3.Range("F13").Formula = "=E13*D13"
$ws3.Range("G13").Value = "Signal generation"
$ws3.Range("H13").Value = "Depending on the laser control interface (Oxxius: SMB). Order 2 per laser line"

# Row 14 - Thorlabs SMA coaxial cable (struck through, not used)
$ws3.Range("A14").Value = "CA2806"
$ws3.Range("B14").Value = "Thorlabs"
$ws3.Range("C14").Value = "SMA Coaxial Cable, SMA Male to BNC Male, 6"" (152 mm)"
$ws3.Range("D14").Value = 1
$ws3.Range("E14").Value = 14
$ws3.Range("F14").Formula = "=E14*D14"
$ws3.Range("G14").Value = "Signal generation"
$ws3.Range("H14").Value = "Depends on laser interface. Oxxius: IO connector must be modified, pin24 (MDL-FSTM) and pin16 (GND Digital), to accept arm switching signal. Ask oxxius to do this when ordering. You want free wire ends. "
$ws3.Range("A14:H14").Font.Strikethrough = $true

# Row 15 - Thorlabs 0.3m BNC cable
$ws3.Range("A15").Value = "2249-C-12"
$ws3.Range("B15").Value = "Thorlabs"
$ws3.Range("C15").Value = "0.3 m BNC Cable"
$ws3.Range("D15").Value = 1
$ws3.Range("E15").Value = 15.25
$ws3.Range("F15").Formula = "=E15*D15"
$ws3.Range("G15").Value = "Signal generation"
$ws3.Range("H15").Value = "necessary for stage control"

# Row 16 - Thorlabs 0.6m BNC cable
$ws3.Range("A16").Value = "2249-C-24"
$ws3.Range("B16").Value = "Thorlabs"
$ws3.Range("C16").Value = "0.6 m BNC Cable"
$ws3.Range("D16").Value = 1
$ws3.Range("E16").Value = 16.5
$ws3.Range("F16").Formula = "=E16*D16"
$ws3.Range("G16").Value = "Signal generation"
$ws3.Range("H16").Value = "Only 1 of these is needed for stage control. Depends on how far the controller is from the terminal block. You could also buy a 	thorlabs CA2024 bnc male to free end cable but this is a bit more expensive… "
$ws3.Range("H16:H18").Merge()
$ws3.Range("H16:H18").HorizontalAlignment = -4108
$ws3.Range("H16:H18").VerticalAlignment = -4108

# Row 17 - Thorlabs 1.5m BNC cable
$ws3.Range("A17").Value = "2249-C-60"
$ws3.Range("B17").Value = "Thorlabs"
$ws3.Range("C17").Value = "1.5 m BNC cable"
$ws3.Range("D17").Value = 1
$ws3.Range("E17").Value = 24
$ws3.Range("F17").Formula = "=E17*D17"
$ws3.Range("G17").Value = "Signal generation"

# Row 18 - Thorlabs 3m BNC cable
$ws3.Range("A18").Value = "2249-C-120"
$ws3.Range("B18").Value = "Thorlabs"
$ws3.Range("C18").Value = "3m BNC cable"
$ws3.Range("D18").Value = 1
$ws3.Range("E18").Value = 24
$ws3.Range("F18").Formula = "=E18*D18"
$ws3.Range("G18").Value = "Signal generation"

# Row 19 - Thorlabs BNC T adapter
$ws3.Range("A19").Value = "T3285"
$ws3.Range("B19").Value = "Thorlabs"
$ws3.Range("C19").Value = "BNC Adapter, T, Female-Male-Female (F-M-F)"
$ws3.Range("D19").Value = 1
$ws3.Range("E19").Value = 11
$ws3.Range("F19").Formula = "=E19*D19"
$ws3.Range("G19").Value = "Signal generation"
$ws3.Range("H19").Value = "Splitting the BNC trigger between Z- and F-cards of ASI controller. Cheaper ones available from other suppliers."
$ws3.Range("A19").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Remove the now-shifted blank spacer row (old row 14, now row 23) that sat
#    between the "Detection" block and the filter-wheel rows.
# ---------------------------------------------------------------------------
$ws3.Rows("23:23").Delete()

# ---------------------------------------------------------------------------
# 4. New "Excitation" section at the bottom of the sheet (rows 26-27), with a
#    struck-through CBLS3P cable-set row.
# ---------------------------------------------------------------------------
$ws3.Range("A26").Value = "Excitation"
$ws3.Range("A26").Font.Name = "Calibri"
$ws3.Range("A26").Font.Size = 11
$ws3.Range("A26").Font.Bold = $true
$ws3.Range("A26").Font.Color = RGB(0, 0, 0)
$ws3.Range("A26").Interior.Color = RGB(242, 242, 242)
$ws3.Range("I26").Value = "Replaces"
$ws3.Range("I26").Font.Name = "Calibri"
$ws3.Range("I26").Font.Size = 11
$ws3.Range("I26").Font.Bold = $true
$ws3.Range("I26").Font.Color = RGB(0, 0, 0)
$ws3.Range("I26").Interior.Color = RGB(242, 242, 242)

$ws3.Range("A27").Value = "CBLS3P"
$ws3.Range("B27").Value = "Thorlabs"
$ws3.Range("C27").Value = "CBLS3P - Command and Power Cables for QS15/20 Scanning Galvanometer Systems and GPS011 Series Power Supply"
$ws3.Range("D27").Value = 2
$ws3.Range("E27").Value = 148
$ws3.Range("F27").Formula = "=E27*D27"
$ws3.Range("G27").Value = "Excitation"

$ws3.Range("A27:G27").Font.Strikethrough = $true

$ws3.Range("H27").Value = "This set of cables was not necessary for me as the power cables came bundled with the GPS011 unit and command cable with free wire ends for connection with the screw terminal block was in the CBLS3F set. CHECK WITH LOCAL THORLABS REP "
$ws3.Range("H27").Font.Name = "Arial"
$ws3.Range("H27").Font.Size = 10
$ws3.Range("H27").Characters(246, 30).Font.Bold = $true

# ---------------------------------------------------------------------------
# 5. Misc view-state bookkeeping to mirror the author's saved selection state.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Main")
$ws1.Range("A136").Select()

$ws2 = $wb.Worksheets.Item("Options")
$ws2.Activate()
$ws2.Range("I6").Select()

$ws3.Activate()
$ws3.Range("A23:XFD23").Select()

Write-Host "done"
